# team roles status update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column F to fit the new, longer status text (stored width = 14)
$ws.Columns.Item(6).ColumnWidth = 13.17

# Update status for the rows that have finished (changed from "-" to "In Document")
$rows = @(6, 12, 15, 18, 24, 27, 30, 33, 36)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "In Document"
}

# Update the active selection to reflect the new cursor position
$ws.Range("H11").Select()
